$wb = $excel.ActiveWorkbook

# ---- Sheet: general ----
$ws = $wb.Worksheets.Item("general")
$ws.Range("B3").Value = 52.48892115189314
$ws.Range("B4").Value = 0.0130000114440918
$ws.Range("B6").Value = 37.36892126633407
$ws.Range("B7").Value = 0
$ws.Range("B8").Value = 0
$ws.Range("B9").Value = 15.11999988555908

# ---- Sheet: alpha ----
$ws = $wb.Worksheets.Item("alpha")
$ws.Range("A3").EntireRow.Delete()
$ws.Range("A2").EntireRow.Delete()

# ---- Sheet: x ----
$ws = $wb.Worksheets.Item("x")
$ws.Range("B2").Value = 2
$ws.Range("B4").Value = 10
$ws.Range("B5").Value = 13
$ws.Range("B6").Value = 9
$ws.Range("B8").Value = 8
$ws.Range("B9").Value = 12
$ws.Range("B10").Value = 4
$ws.Range("B11").Value = 1
$ws.Range("B12").Value = 11
$ws.Range("B13").Value = 5
$ws.Range("B14").Value = 7

# ---- Sheet: TBar ----
$ws = $wb.Worksheets.Item("TBar")
$ws.Range("B3").Value = 32.61192465059682
$ws.Range("B4").Value = 30
$ws.Range("B6").Value = 35.60566989569675
$ws.Range("B7").Value = 36.09699127290008
$ws.Range("B9").Value = 32.01159140980468
$ws.Range("B10").Value = 37.91791988218178
$ws.Range("B11").Value = 34.14711948224307
$ws.Range("B12").Value = 30
$ws.Range("B13").Value = 43.13422125384351
$ws.Range("B14").Value = 37.28210730097497
$ws.Range("B15").Value = 37.94859027624736

# ---- Sheet: y ----
$ws = $wb.Worksheets.Item("y")
$ws.Range("A3").EntireRow.Delete()
$ws.Range("A2").EntireRow.Delete()

# ---- Sheet: Q ----
$ws = $wb.Worksheets.Item("Q")
$ws.Range("C7").Value = 250.9700000000009
$ws.Range("C8").Value = 260.9900000000009
$ws.Range("C9").Value = 252.975000000001
$ws.Range("C10").Value = 269.580000000001
$ws.Range("C11").Value = 250.575000000001
$ws.Range("C12").Value = 81.47500000000073
$ws.Range("C13").Value = 80.68000000000073
$ws.Range("C14").Value = 84.71500000000073
$ws.Range("C15").Value = 80.43500000000073
$ws.Range("C16").Value = 87.34500000000072
$ws.Range("C21").Value = 39.43499999999942
$ws.Range("C22").Value = 226.9299999999988
$ws.Range("C23").Value = 228.4
$ws.Range("C24").Value = 211.0249999999988
$ws.Range("C25").Value = 230.2799999999988
$ws.Range("C26").Value = 215.0299999999988
$ws.Range("C27").Value = 236.8400000000015
$ws.Range("C28").Value = 244.9650000000015
$ws.Range("C29").Value = 236.8400000000015
$ws.Range("C30").Value = 243.8100000000014
$ws.Range("C31").Value = 244.9650000000015
$ws.Range("C37").Value = 141.0250000000001
$ws.Range("C38").Value = 143.4
$ws.Range("C39").Value = 139.7050000000002
$ws.Range("C40").Value = 150.4250000000002
$ws.Range("C41").Value = 134.7700000000002
$ws.Range("C42").Value = 294.8549999999982
$ws.Range("C43").Value = 307.5599999999982
$ws.Range("C44").Value = 270.8449999999982
$ws.Range("C45").Value = 294.05
$ws.Range("C46").Value = 273.9599999999982
$ws.Range("C47").Value = 166.9600000000012
$ws.Range("C48").Value = 168.6450000000012
$ws.Range("C49").Value = 164.4300000000013
$ws.Range("C50").Value = 171.1650000000012
$ws.Range("C51").Value = 172.0750000000012
$ws.Range("C52").Value = 57.95
$ws.Range("C53").Value = 58.67999999999927
$ws.Range("C54").Value = 61.72999999999927
$ws.Range("C55").Value = 60.65499999999928
$ws.Range("C56").Value = 52.91499999999927
$ws.Range("C57").Value = 294.8549999999982
$ws.Range("C58").Value = 307.5599999999982
$ws.Range("C59").Value = 270.8449999999982
$ws.Range("C60").Value = 294.05
$ws.Range("C61").Value = 273.9599999999982
$ws.Range("C62").Value = 236.8400000000015
$ws.Range("C63").Value = 244.9650000000015
$ws.Range("C64").Value = 236.8400000000015
$ws.Range("C65").Value = 243.8100000000014
$ws.Range("C66").Value = 244.9650000000015
$ws.Range("C67").Value = 250.9700000000009

# ---- Sheet: R ----
$ws = $wb.Worksheets.Item("R")
$ws.Range("C3").Value = 7.559999942779541
$ws.Range("C8").Value = 0
$ws.Range("C10").Value = 0

# ---- Sheet: rho ----
$ws = $wb.Worksheets.Item("rho")
$ws.Range("A2").Value = 8
$ws.Range("A3").EntireRow.Delete()
